$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B and label its header "Loc".
# This shifts the existing "Tag no", "Problem found", "Action Taken",
# "complain", "date" headers one column to the right (B->C, C->D, D->E, E->F, F->G).
$ws.Range("B1").EntireColumn.Insert()

$ws.Range("B2").Value = "Loc"

# Column width tweaks that came along with the paste (now landing on the
# shifted "Problem found" / "Action Taken" columns).
$ws.Range("D1").EntireColumn.ColumnWidth = 15.85546875
$ws.Range("E1").EntireColumn.ColumnWidth = 37.85546875

# Reflect the pasted-range selection left behind on the sheet.
$ws.Range("B3:G52").Select()
$excel.ActiveWindow.RangeSelection.Activate()

# Page setup as captured by the diff.
$ws.PageSetup.Orientation = 1
